# Add new rows (robustness-cycles data) to both the NBR and BAR sheets.
# Row 2 already exists with cutoff=0; new rows continue cutoff=1..15 with
# NumCycles 6..20, mirroring the existing format (col A bold/bordered style,
# cols B/C unformatted).

$wb = $excel.ActiveWorkbook

$nbrValues = @(
    @(853),
    @(818),
    @(821),
    @(821),
    @(847),
    @(847),
    @(843),
    @(848),
    @(843),
    @(843),
    @(840),
    @(820),
    @(817),
    @(809),
    @(802)
)

$barValues = @(
    @(973),
    @(998),
    @(1000),
    @(999),
    @(971),
    @(962),
    @(959),
    @(931),
    @(929),
    @(919),
    @(924),
    @(917),
    @(917),
    @(916),
    @(913)
)

function Fill-Sheet($ws, $values) {
    # Copy the formatting of the existing data row (row 2) down to the
    # new rows so column A keeps its bold/bordered style.
    $lastRow = 2 + $values.Length
    $ws.Range("A2:C2").Copy()
    $ws.Range("A3:C$lastRow").PasteSpecial(-4122) # xlPasteFormats

    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = 3 + $i
        $ws.Cells.Item($row, 1).Value = $i + 1
        $ws.Cells.Item($row, 2).Value = $i + 6
        $ws.Cells.Item($row, 3).Value = $values[$i][0]
    }
}

$wsNbr = $wb.Worksheets.Item("NBR")
Fill-Sheet $wsNbr $nbrValues

$wsBar = $wb.Worksheets.Item("BAR")
Fill-Sheet $wsBar $barValues
